$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: "Scalpel Accuracy:" text moves from C34 to E34, and the 100 value moves from D34 to F34.
$ws.Range("C34").Value = $null
$ws.Range("D34").Value = $null
$ws.Range("E34").Value = "Scalpel Accuracy:"
$ws.Range("F34").Value = 100

# Row 35: fix the label text
$ws.Range("E35").Value = "Accuracy vs PyType"
